# Sheet order in this workbook:
#   1 Funciones_Objetivo          (unchanged)
#   2 Restricciones_del_lider     (unchanged)
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha
# NOTE: sheet names "Vector_bf" / "Vector_BF" differ only by case, and
# Worksheets.Item(name) lookups here are case-insensitive, so sheets are
# addressed by their (unambiguous) 1-based index instead of by name.

$wb = $excel.ActiveWorkbook

# ---- Restricciones_del_follower (index 3) ----
$ws = $wb.Worksheets.Item(3)

# Force text storage for the cells we touch (these values look numeric but
# must remain stored as text/shared-strings, matching the source data).
# (E4's value does not actually change - "0" stays "0" - so it is left alone
# to avoid touching its style/format for no reason.)
$ws.Range("A2:B6").NumberFormat = "@"
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("E2:E3").NumberFormat = "@"
$ws.Range("E5:E6").NumberFormat = "@"
$ws.Range("F2:F6").NumberFormat = "@"

$ws.Range("A2").Value = "-2.666551724137933 - x + 1.2068965517241381y"
$ws.Range("B2").Value = "4.666551724137933"
$ws.Range("D2").Value = "0.93"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "2.1"

$ws.Range("A3").Value = "-19.01706896551724 - 0.25x + 3.275862068965517y"
$ws.Range("B3").Value = "17.01706896551724"
$ws.Range("D3").Value = "0.41"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "5.699999999999999"

$ws.Range("A4").Value = "-4.78 + x"
$ws.Range("B4").Value = "-3.2199999999999998"
$ws.Range("D4").Value = "0.7"
$ws.Range("F4").Value = "10.0"

$ws.Range("A5").Value = "-16.69172413793104 + x + 1.8965517241379315y"
$ws.Range("B5").Value = "14.481724137931039"
$ws.Range("D5").Value = "0.36"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "3.3000000000000003"

$ws.Range("A6").Value = "-0.17"
$ws.Range("B6").Value = "0"
$ws.Range("D6").Value = "0.79"
$ws.Range("E6").Value = "8.299999999999999"
$ws.Range("F6").Value = "0"

# ---- Punto_modificado (index 4) - modified point (x, y) ----
$ws2 = $wb.Worksheets.Item(4)
$ws2.Range("A2:B2").NumberFormat = "@"
$ws2.Range("A2").Value = "4.78"
$ws2.Range("B2").Value = "6.17"

# ---- Vector_bf (index 5) ----
$ws3 = $wb.Worksheets.Item(5)
$ws3.Range("A2").NumberFormat = "@"
$ws3.Range("A2").Value = "-2.148275862068966"

# ---- Vector_BF (index 6) ----
$ws4 = $wb.Worksheets.Item(6)
$ws4.Range("A2:A3").NumberFormat = "@"
$ws4.Range("A2").Value = "-1.0"
$ws4.Range("A3").Value = "-1.0"

# ---- Vector_Alpha (index 7) - A2 is a genuine numeric cell (not text) ----
$ws5 = $wb.Worksheets.Item(7)
$ws5.Range("A2").Value = 1.7399999999999998
